$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 5 (pushes existing rows 5-33 down to 6-34)
$ws.Rows.Item(5).Insert()

# Populate the new Issue_033 row
$ws.Range("A5").Value = "Issue_033"
$ws.Range("B5").Value = "Rearchitect Build Evaluation String to use the test class"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "Open"
$ws.Range("E5").Value = "\Plugins\Utility\Build Evaluation String.vi"

# Copy formatting (borders/wrap/styles) from the row below, which held the old row 5 data
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-apply values since paste-special formats may have touched formulas/values unexpectedly
$ws.Range("A5").Value = "Issue_033"
$ws.Range("B5").Value = "Rearchitect Build Evaluation String to use the test class"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "Open"
$ws.Range("E5").Value = "\Plugins\Utility\Build Evaluation String.vi"
$ws.Range("F5").Value = ""

# Update the view/selection to match the final state
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("B32").Select()
